# Commit: "add new results to be analysed"
#
# The simulation produced 20 more samples (a new obstacle configuration),
# so 20 more columns (V:AO) of results are appended next to the existing
# B:U block on Sheet1:
#   - row 2 keeps storing the sampled angles (0 .. 2*pi), the same 20
#     values used for B2:U2 are written again for V2:AO2;
#   - row 3 stores the new pass/fail (1/0) outcomes for those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: angle samples (0 .. 2*pi), same pattern as B2:U2
$row2Values = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# Row 3: new pass/fail (1/0) results for columns V:AO
$row3Values = @(1, 1, 0, 1, 1, 1, 0, 0, 0, 1, 1, 1, 0, 0, 0, 0, 0, 0, 1, 1)

# Columns V (22) through AO (41)
$startCol = 22
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
    $ws.Cells.Item(3, $col).Value = $row3Values[$i]
}

# Scroll/select so the newly added columns are in view, mirroring the
# author's on-screen state after pasting the new results
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.Left = 1065
$win.Top = 3120
$ws.Range("Y10").Select()

$wb.Save()
